# Adding new estimates (pop; alone or in combo)
#
# The lookup table's "pop" row (row 2) originally represented the combined
# "Alone or in combination" estimate but was just labeled "US residents,
# 2016". We now:
#   1. Relabel row 2 ("pop") to be explicit: "US residents, 2016 (Alone or
#      in combo)".
#   2. Insert a brand-new row right after it for the new "pop_v2" variable,
#      labeled "US residents, 2016 (Alone)", still under the "Total
#      Population" topic group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes everything from the old row 3 down
# to row 4, carrying formatting/styles along with it).
$ws.Rows.Item(3).EntireRow.Insert() | Out-Null

# Populate the newly inserted row 3 with the new "pop_v2" variable.
$ws.Range("A3").Value = "pop_v2"
$ws.Range("B3").Value = "US residents, 2016 (Alone)"
$ws.Range("C3").Value = "Total Population"

# Update row 2's label to clarify it is the "alone or in combo" estimate.
$ws.Range("B2").Value = "US residents, 2016 (Alone or in combo)"

# Match the author's final cursor position.
$ws.Range("B3").Select() | Out-Null
